$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 67.8679804978435
$ws.Range("C2").Value = 66.47906002391
$ws.Range("D2").Value = 69.256900971777
$ws.Range("C3").Value = 63.3414374246431
$ws.Range("D3").Value = 68.0972229127337
$ws.Range("B11").Value = 67.987520077397
$ws.Range("C11").Value = 64.4000884639409
$ws.Range("D11").Value = 71.5749516908532
$ws.Range("B12").Value = 74.0846263450239
$ws.Range("C12").Value = 70.2362618026456
$ws.Range("D12").Value = 77.9329908874022
$ws.Range("B14").Value = 64.0735542026883
$ws.Range("C14").Value = 61.0923553714412
$ws.Range("D14").Value = 67.0547530339354
$ws.Range("B15").Value = 55.2953815419129
$ws.Range("C15").Value = 53.8133973581825
$ws.Range("D15").Value = 56.7773657256433
$ws.Range("C16").Value = 58.6004210267478
$ws.Range("D16").Value = 63.5825322998283
$ws.Range("B24").Value = 56.8474042011213
$ws.Range("C24").Value = 53.0475549407763
$ws.Range("D24").Value = 60.6472534614663
$ws.Range("B25").Value = 48.4302557813727
$ws.Range("C25").Value = 44.0801483613297
$ws.Range("D25").Value = 52.7803632014157
$ws.Range("B27").Value = 48.666425997025
$ws.Range("C27").Value = 45.5560116954381
$ws.Range("D27").Value = 51.7768402986118
$ws.Range("B28").Value = 64.7240480131028
$ws.Range("C28").Value = 63.0045795926267
$ws.Range("D28").Value = 66.4435164335789
$ws.Range("C29").Value = 66.154032049576
$ws.Range("D29").Value = 71.4780870671454
$ws.Range("B37").Value = 66.5933688883394
$ws.Range("C37").Value = 62.4403846149291
$ws.Range("D37").Value = 70.7463531617498
$ws.Range("B38").Value = 57.2690723237937
$ws.Range("C38").Value = 50.7470363646041
$ws.Range("D38").Value = 63.7911082829833
$ws.Range("B40").Value = 60.3471262597791
$ws.Range("C40").Value = 56.8080642142266
$ws.Range("D40").Value = 63.8861883053316
$ws.Range("B41").Value = 67.2054760884641
$ws.Range("C41").Value = 65.5001731723724
$ws.Range("D41").Value = 68.9107790045557
$ws.Range("B42").Value = 73.150166637206
$ws.Range("C42").Value = 70.5221212688109
$ws.Range("D42").Value = 75.7782120056012
$ws.Range("B50").Value = 67.292535125419
$ws.Range("C50").Value = 63.1364588496784
$ws.Range("D50").Value = 71.4486114011597
$ws.Range("B51").Value = 63.9033952437169
$ws.Range("C51").Value = 58.2061893987969
$ws.Range("D51").Value = 69.6006010886369
$ws.Range("B53").Value = 58.9468339412676
$ws.Range("C53").Value = 55.254083851168
$ws.Range("D53").Value = 62.6395840313673
$ws.Range("B54").Value = 68.518253119246
$ws.Range("C54").Value = 66.8882724435162
$ws.Range("D54").Value = 70.1482337949759
$ws.Range("C55").Value = 70.7689055039649
$ws.Range("D55").Value = 75.9022703411351
$ws.Range("B63").Value = 70.8149579691965
$ws.Range("C63").Value = 66.8120913855674
$ws.Range("D63").Value = 74.8178245528256
$ws.Range("B64").Value = 60.4866837966034
$ws.Range("C64").Value = 54.4980227132752
$ws.Range("D64").Value = 66.4753448799316
$ws.Range("B66").Value = 63.6828257869627
$ws.Range("C66").Value = 60.3796020159564
$ws.Range("D66").Value = 66.986049557969
$ws.Range("B67").Value = 61.895417279981
$ws.Range("C67").Value = 60.1645023547081
$ws.Range("D67").Value = 63.6263322052539
$ws.Range("C68").Value = 65.1898278582191
$ws.Range("D68").Value = 70.736821035104
$ws.Range("B76").Value = 63.3861854984241
$ws.Range("C76").Value = 59.2235413127094
$ws.Range("D76").Value = 67.5488296841387
$ws.Range("B77").Value = 54.11746546684
$ws.Range("C77").Value = 48.0133810262521
$ws.Range("D77").Value = 60.2215499074278
$ws.Range("B79").Value = 56.0767240057917
$ws.Range("C79").Value = 52.589586582689
$ws.Range("D79").Value = 59.5638614288945